# Update marksheet correction counts / totals
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking row: right-answer mark value 3 -> 5
$ws.Range("B11").Value = 5

# Total row: total marks 63 -> 105, and corr/total text 58/84 -> 105/140
$ws.Range("B12").Value = 105
$ws.Range("E12").Value = "105/140"
